$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a brand-new column before AI. Everything that used to live in
#    columns AI:AY (illinois .. ontario) shifts one column right to AJ:AZ.
# ---------------------------------------------------------------------------
$ws.Columns("AI").Insert()

# ---------------------------------------------------------------------------
# 2) Header row: new country columns "cambodia" (now AI1) and "arizona"
#    (appended at the very end, BA1). BA is beyond the sheet's declared
#    <cols> range, so clone AZ1's format first to keep the same header style.
# ---------------------------------------------------------------------------
$ws.Range("AZ1").Copy()
$ws.Range("BA1").PasteSpecial(-4122)

$ws.Range("AI1").Value = "cambodia"
$ws.Range("BA1").Value = "arizona"

# ---------------------------------------------------------------------------
# 3) New data row (row 10) for 2020-01-27. Clone row 9's formatting first
#    (so cells keep the same text-number style already used by the table)
#    and then fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A9:AU9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("A10").Value = 43857

$ws.Range("B10").Value = "70"
$ws.Range("C10").Value = "72-0-2"
$ws.Range("D10").Value = "110-0-0-1"
$ws.Range("E10").Value = "56"
$ws.Range("F10").Value = "14"
$ws.Range("G10").Value = "151-0-3"
$ws.Range("H10").Value = "46"
$ws.Range("I10").Value = "7"
$ws.Range("J10").Value = "33-0-0-1"
$ws.Range("K10").Value = "18-0-0-1"
$ws.Range("L10").Value = "21-0-0-1"
$ws.Range("M10").Value = "128"
$ws.Range("N10").Value = "8"
$ws.Range("O10").Value = "1423-0-45-76"
$ws.Range("P10").Value = "100"
$ws.Range("Q10").Value = "11"
$ws.Range("R10").Value = "47-0-1"
$ws.Range("S10").Value = "48-0-1"
$ws.Range("T10").Value = "6"
$ws.Range("U10").Value = "27"
$ws.Range("V10").Value = "6"
$ws.Range("W10").Value = "7"
$ws.Range("X10").Value = "4"
$ws.Range("Y10").Value = "35"
$ws.Range("Z10").Value = "75"
$ws.Range("AA10").Value = "53-0-3-1"
$ws.Range("AB10").Value = "13"
$ws.Range("AC10").Value = "69"
$ws.Range("AD10").Value = "5"
$ws.Range("AE10").Value = "22"
$ws.Range("AF10").Value = "5"
$ws.Range("AG10").Value = "26"
$ws.Range("AH10").Value = "128-0-1"
$ws.Range("AI10").Value = "1"
$ws.Range("AJ10").Value = "1"
$ws.Range("AK10").Value = "1"
$ws.Range("AL10").Value = "2"
$ws.Range("AM10").Value = "4-0-1"
$ws.Range("AN10").Value = "8-0-2"
$ws.Range("AO10").Value = "4"
$ws.Range("AP10").Value = "4"
$ws.Range("AQ10").Value = "2"
$ws.Range("AR10").Value = "3"
$ws.Range("AS10").Value = "5"
$ws.Range("AT10").Value = "1"
$ws.Range("AU10").Value = "4"

$ws.Range("BA10").Value = 1
